$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.329689860343933
$ws.Range("B1").Value = 1.600885152816772
$ws.Range("C1").Value = 2.157634258270264
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.175687313079834
